$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-typed cell updates (avoid Excel auto-converting numeric-looking
# strings to floating point values / percentages) by using a leading apostrophe
# and then resetting the style so no extra "Text" number-format style sticks.
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "43.042.16"
Set-TextCell $ws "E2" "  -6.10%  "

Set-TextCell $ws "D3" "2.549.96"
Set-TextCell $ws "E3" "  -2.12%  "

Set-TextCell $ws "E4" "  +0.09%  "

Set-TextCell $ws "D5" "299.40"
Set-TextCell $ws "E5" "  -3.09%  "

Set-TextCell $ws "D6" "92.54"
Set-TextCell $ws "E6" "  -6.28%  "

Set-TextCell $ws "E7" "  -3.36%  "

Set-TextCell $ws "D9" "0.549"
Set-TextCell $ws "E9" "  -4.73%  "

Set-TextCell $ws "E10" "  -7.28%  "

Set-TextCell $ws "D11" "0.0805"
Set-TextCell $ws "E11" "  -3.96%  "

Set-TextCell $ws "D12" "7.65"
Set-TextCell $ws "E12" "  -4.84%  "

Set-TextCell $ws "D13" "0.112"
Set-TextCell $ws "E13" "  +4.77%  "

Set-TextCell $ws "D14" "2.939.24"
Set-TextCell $ws "E14" "  -1.94%  "

Set-TextCell $ws "D15" "2.544.10"
Set-TextCell $ws "E15" "  -2.15%  "

Set-TextCell $ws "D16" "0.869"
Set-TextCell $ws "E16" "  -4.44%  "

Set-TextCell $ws "D17" "14.10"
Set-TextCell $ws "E17" "  -4.51%  "

Set-TextCell $ws "D18" "43.061.34"
Set-TextCell $ws "E18" "  -6.04%  "

Set-TextCell $ws "D19" "13.05"
Set-TextCell $ws "E19" "  +3.14%  "

Set-TextCell $ws "D20" "0.0₃0982"
Set-TextCell $ws "E20" "  -2.77%  "

Set-TextCell $ws "E21" "  -1.28%  "

Set-TextCell $ws "D22" "71.73"
Set-TextCell $ws "E22" "  -2.83%  "

Set-TextCell $ws "D23" "256.44"
Set-TextCell $ws "E23" "  -9.88%  "

Set-TextCell $ws "D24" "2.90"
Set-TextCell $ws "E24" "  -3.89%  "

Set-TextCell $ws "D25" "2.12"
Set-TextCell $ws "E25" "  -5.57%  "

Set-TextCell $ws "D26" "29.11"
Set-TextCell $ws "E26" "  -0.35%  "

Set-TextCell $ws "E27" "  -0.01%  "

Set-TextCell $ws "D28" "10.03"
Set-TextCell $ws "E28" "  -4.83%  "

Set-TextCell $ws "D29" "37.29"
Set-TextCell $ws "E29" "  -3.13%  "

Set-TextCell $ws "E30" "  -5.75%  "

Set-TextCell $ws "D31" "5.94"
Set-TextCell $ws "E31" "  -4.30%  "

Set-TextCell $ws "D32" "153.09"
Set-TextCell $ws "E32" "  -2.62%  "

Set-TextCell $ws "D33" "2.16"
Set-TextCell $ws "E33" "  -4.47%  "

Set-TextCell $ws "E34" "  -1.56%  "

Set-TextCell $ws "E35" "  -7.12%  "

Set-TextCell $ws "D36" "0.0800"
Set-TextCell $ws "E36" "  -3.72%  "

Set-TextCell $ws "E37" "  -5.15%  "

Set-TextCell $ws "E38" "  -2.15%  "

Set-TextCell $ws "D39" "17.04"
Set-TextCell $ws "E39" "  +7.22%  "

Set-TextCell $ws "D40" "23.26"
Set-TextCell $ws "E40" "  +8.61%  "

Set-TextCell $ws "D41" "3.43"
Set-TextCell $ws "E41" "  -2.47%  "

Set-TextCell $ws "E42" "  -2.62%  "

Set-TextCell $ws "E43" "  -4.32%  "

Set-TextCell $ws "D44" "2.076.77"
Set-TextCell $ws "E44" "  -1.25%  "

Set-TextCell $ws "D45" "0.999"
Set-TextCell $ws "E45" "  +0.05%  "

Set-TextCell $ws "D46" "84.66"
Set-TextCell $ws "E46" "  -9.74%  "

Set-TextCell $ws "E47" "  -3.03%  "

Set-TextCell $ws "E48" "  +1.79%  "

Set-TextCell $ws "D49" "2.796.28"
Set-TextCell $ws "E49" "  -1.83%  "

Set-TextCell $ws "D50" "104.82"
Set-TextCell $ws "E50" "  -3.64%  "

Set-TextCell $ws "D51" "1.67"
Set-TextCell $ws "E51" "  -3.71%  "

